$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1, Q1 - values plus same formatting as the rest of row 1 (style s="1")
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update existing columns I, K, M, O for rows 2-25, and add new columns P, Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # column I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # column K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # column M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # column O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # column P: new, value 2
    $ws.Cells.Item($r, 17).Value = 2   # column Q: new, value 2
}
